$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 32
$ws.Range("H32").Value = 1279
$ws.Range("I32").Value = 370
$ws.Range("J32").Value = 1733.5
$ws.Range("K32").Value = 370
$ws.Range("L32").Value = 1733.5
$ws.Range("M32").Value = -44
$ws.Range("N32").Value = -2385.5

# ALC row 55
$ws.Range("H55").Value = 582.1667
$ws.Range("I55").Value = 668.6
$ws.Range("J55").Value = 150
$ws.Range("K55").Value = 668.6
$ws.Range("L55").Value = 150
$ws.Range("M55").Value = -454.6
$ws.Range("N55").Value = -578

# ALC row 64
$ws.Range("H64").Value = 3906.848
$ws.Range("I64").Value = 3915.8333
$ws.Range("J64").Value = 3874.5
$ws.Range("K64").Value = 3915.8333
$ws.Range("L64").Value = 3874.5
$ws.Range("M64").Value = -3667.8333
$ws.Range("N64").Value = -4370.5

# ALC row 67
$ws.Range("H67").Value = 3906.848
$ws.Range("I67").Value = 3915.8333
$ws.Range("J67").Value = 3874.5
$ws.Range("K67").Value = 3915.8333
$ws.Range("L67").Value = 3874.5
$ws.Range("M67").Value = -3057.8333
$ws.Range("N67").Value = -5590.5

# ALC row 103
$ws.Range("H103").Value = 650666
$ws.Range("I103").Value = 485
$ws.Range("J103").Value = 813211.25
$ws.Range("K103").Value = 1455
$ws.Range("L103").Value = 2439633.75
$ws.Range("M103").Value = -869
$ws.Range("N103").Value = -2440805.75

# ALC row 112
$ws.Range("H112").Value = 868.17145
$ws.Range("J112").Value = 901.7778
$ws.Range("L112").Value = 2705.3334
$ws.Range("N112").Value = -4921.3334

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 4027.3948
$ws.Range("I32").Value = 3333.5483
$ws.Range("J32").Value = 7100.143
$ws.Range("K32").Value = 3333.5483
$ws.Range("L32").Value = 7100.143
$ws.Range("M32").Value = -3046.5483
$ws.Range("N32").Value = -7674.143

# ARM row 45
$ws.Range("H45").Value = 8892.526
$ws.Range("I45").Value = 11502.429
$ws.Range("J45").Value = 1584.8
$ws.Range("K45").Value = 11502.429
$ws.Range("L45").Value = 1584.8
$ws.Range("M45").Value = -11125.429
$ws.Range("N45").Value = -2338.8

# ARM row 132
$ws.Range("H132").Value = 2948.1072
$ws.Range("I132").Value = 2459.8108
$ws.Range("J132").Value = 3899
$ws.Range("K132").Value = 7379.432400000001
$ws.Range("L132").Value = 11697
$ws.Range("M132").Value = -4849.432400000001
$ws.Range("N132").Value = -16757

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Range("H31").Value = 2997.5386
$ws.Range("I31").Value = 1556.1052
$ws.Range("J31").Value = 3827.4546
$ws.Range("K31").Value = 1556.1052
$ws.Range("L31").Value = 3827.4546
$ws.Range("M31").Value = -1261.1052
$ws.Range("N31").Value = -4417.4546

# CRP row 34
$ws.Range("H34").Value = 2997.5386
$ws.Range("I34").Value = 1556.1052
$ws.Range("J34").Value = 3827.4546
$ws.Range("K34").Value = 1556.1052
$ws.Range("L34").Value = 3827.4546
$ws.Range("M34").Value = -1354.1052
$ws.Range("N34").Value = -4231.4546

$ws = $wb.Worksheets.Item("CUL")
# CUL row 5
$ws.Range("H5").Value = 9925.546
$ws.Range("I5").Value = 50392
$ws.Range("J5").Value = 933
$ws.Range("K5").Value = 151176
$ws.Range("L5").Value = 2799
$ws.Range("M5").Value = -151064
$ws.Range("N5").Value = -3023

# CUL row 131
$ws.Range("H131").Value = 17743104
$ws.Range("I131").Value = 10000442
$ws.Range("J131").Value = 19232076
$ws.Range("K131").Value = 30001326
$ws.Range("L131").Value = 57696228
$ws.Range("M131").Value = -29996286
$ws.Range("N131").Value = -57706308

# CUL row 135
$ws.Range("H135").Value = 9925.546
$ws.Range("I135").Value = 50392
$ws.Range("J135").Value = 933
$ws.Range("K135").Value = 453528
$ws.Range("L135").Value = 8397
$ws.Range("M135").Value = -450993
$ws.Range("N135").Value = -13467

$ws = $wb.Worksheets.Item("GSM")
# GSM row 38
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

# GSM row 57
$ws.Range("H57").Value = 16133.5
$ws.Range("J57").Value = 16133.5
$ws.Range("L57").Value = 16133.5
$ws.Range("N57").Value = -17773.5

# GSM row 122
$ws.Range("H122").Value = 36719344
$ws.Range("I122").Value = 53242440
$ws.Range("J122").Value = 1360.1111
$ws.Range("K122").Value = 159727320
$ws.Range("L122").Value = 4080.3333
$ws.Range("M122").Value = -159724870
$ws.Range("N122").Value = -8980.3333

# GSM row 132
$ws.Range("H132").Value = 3928.075
$ws.Range("I132").Value = 4711.2
$ws.Range("J132").Value = 3667.0334
$ws.Range("K132").Value = 14133.6
$ws.Range("L132").Value = 11001.1002
$ws.Range("M132").Value = -11603.6
$ws.Range("N132").Value = -16061.1002

$ws = $wb.Worksheets.Item("LTW")
# LTW row 7
$ws.Range("H7").Value = 58082.945
$ws.Range("I7").Value = 93053.45
$ws.Range("J7").Value = 3129.2856
$ws.Range("K7").Value = 93053.45
$ws.Range("L7").Value = 3129.2856
$ws.Range("M7").Value = -92941.45
$ws.Range("N7").Value = -3353.2856

# LTW row 22
$ws.Range("H22").Value = 6902778
$ws.Range("I22").Value = 18730968
$ws.Range("K22").Value = 18730968
$ws.Range("M22").Value = -18730673

# LTW row 27
$ws.Range("H27").Value = 6902778
$ws.Range("I27").Value = 18730968
$ws.Range("K27").Value = 18730968
$ws.Range("M27").Value = -18730861

# LTW row 40
$ws.Range("H40").Value = 83336210
$ws.Range("I40").Value = 100002650
$ws.Range("K40").Value = 100002650
$ws.Range("M40").Value = -100002514

# LTW row 46
$ws.Range("H46").Value = 19609210
$ws.Range("I46").Value = 33334298
$ws.Range("J46").Value = 1942.7142
$ws.Range("K46").Value = 33334298
$ws.Range("L46").Value = 1942.7142
$ws.Range("M46").Value = -33334110
$ws.Range("N46").Value = -2318.7142

# LTW row 126
$ws.Range("H126").Value = 58082.945
$ws.Range("I126").Value = 93053.45
$ws.Range("J126").Value = 3129.2856
$ws.Range("K126").Value = 279160.35
$ws.Range("L126").Value = 9387.856800000001
$ws.Range("M126").Value = -276690.35
$ws.Range("N126").Value = -14327.8568

# LTW row 127
$ws.Range("H127").Value = 34998
$ws.Range("J127").Value = 34998
$ws.Range("L127").Value = 34998
$ws.Range("N127").Value = -44918

# LTW row 136
$ws.Range("H136").Value = 3914.394
$ws.Range("I136").Value = 2467.318
$ws.Range("J136").Value = 6808.5454
$ws.Range("K136").Value = 7401.954000000001
$ws.Range("L136").Value = 20425.6362
$ws.Range("M136").Value = -4851.954000000001
$ws.Range("N136").Value = -25525.6362

$ws = $wb.Worksheets.Item("WVR")
# WVR row 81
$ws.Range("H81").Value = 1910.3
$ws.Range("I81").Value = 1633
$ws.Range("K81").Value = 3266
$ws.Range("M81").Value = -2205

# WVR row 84
$ws.Range("H84").Value = 1910.3
$ws.Range("I84").Value = 1633
$ws.Range("K84").Value = 16330
$ws.Range("M84").Value = -11026

# WVR row 122
$ws.Range("H122").Value = 3134.875
$ws.Range("I122").Value = 3020.4167
$ws.Range("J122").Value = 3249.3333
$ws.Range("K122").Value = 9061.250100000001
$ws.Range("L122").Value = 9747.999899999999
$ws.Range("M122").Value = -6611.250100000001
$ws.Range("N122").Value = -14647.9999

# WVR row 132
$ws.Range("H132").Value = 1983.8
$ws.Range("I132").Value = 1808.1482
$ws.Range("J132").Value = 2576.625
$ws.Range("K132").Value = 5424.444600000001
$ws.Range("L132").Value = 7729.875
$ws.Range("M132").Value = -2894.444600000001
$ws.Range("N132").Value = -12789.875

Write-Host "Applied all Ixion_Profits.xlsx updates"
